$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 'fn foo() -> usize {
   let mut a = Vec::new();
   a.len()
}'
$ws.Range("F2").Value = 'fn foo() -> usize {
    let mut a = Vec::new();
    a.len()
}'
$ws.Range("E5").Value = 'asm!("lea {}, [{}]", lateout(reg) _, in(reg) ptr);'
$ws.Range("F5").Value = 'asm!("lea ({}), {}", in(reg) ptr, lateout(reg) _, options(att_syntax));'
$ws.Range("E6").Value = 'asm!("lea ({}), {}", in(reg) ptr, lateout(reg) _, options(att_syntax));'
$ws.Range("F6").Value = 'asm!("lea {}, [{}]", lateout(reg) _, in(reg) ptr);'
$ws.Range("E8").Value = 'assert!(r.is_ok());
assert!(r.is_err());'
$ws.Range("E11").Value = 'fn not_quite_hot_code(..) { ... }'
$ws.Range("E12").Value = 'extern crate foo;
use foo::bar;'
$ws.Range("F12").Value = 'use foo::baz;
extern crate baz;'
$ws.Range("E13").Value = 'fn something_else() { /* ... */ }'
$ws.Range("E14").Value = '
fn not_quite_good_code() { }'
$ws.Range("F14").Value = '// Good (as inner attribute)

fn this_is_fine() { }

// or

// Good (as outer attribute)
fn this_is_fine_too() { }'
$ws.Range("E16").Value = ''
$ws.Range("F16").Value = ''
$ws.Range("E17").Value = 'fn main() { }'
$ws.Range("F17").Value = 'fn main() { }'
$ws.Range("E18").Value = 'fn conditional() { }'
$ws.Range("F18").Value = 'fn conditional() { }

// or

fn conditional() { }

Check the Rust Reference for more details.'
$ws.Range("E19").Value = ''
$ws.Range("F19").Value = ''
$ws.Range("E20").Value = 'pub struct Bar;'
$ws.Range("F20").Value = 'pub struct Bar;'
$ws.Range("E21").Value = 'async fn foo(x: &Mutex<u32>) {
  let mut guard = x.lock().unwrap();
  *guard += 1;
  baz().await;
}

async fn bar(x: &Mutex<u32>) {
  let mut guard = x.lock().unwrap();
  *guard += 1;
  drop(guard); // explicit drop
  baz().await;
}'
$ws.Range("F21").Value = 'async fn foo(x: &Mutex<u32>) {
  {
    let mut guard = x.lock().unwrap();
    *guard += 1;
  }
  baz().await;
}

async fn bar(x: &Mutex<u32>) {
  {
    let mut guard = x.lock().unwrap();
    *guard += 1;
  } // guard dropped here at end of scope
  baz().await;
}'
$ws.Range("E22").Value = 'async fn foo(x: &RefCell<u32>) {
  let mut y = x.borrow_mut();
  *y += 1;
  baz().await;
}

async fn bar(x: &RefCell<u32>) {
  let mut y = x.borrow_mut();
  *y += 1;
  drop(y); // explicit drop
  baz().await;
}'
$ws.Range("F22").Value = 'async fn foo(x: &RefCell<u32>) {
  {
     let mut y = x.borrow_mut();
     *y += 1;
  }
  baz().await;
}

async fn bar(x: &RefCell<u32>) {
  {
    let mut y = x.borrow_mut();
    *y += 1;
  } // y dropped here at end of scope
  baz().await;
}'
$ws.Range("E23").Value = 'await-holding-invalid-types = [
  # You can specify a type name
  "CustomLockType",
  # You can (optionally) specify a reason
  { path = "OtherCustomLockType", reason = "Relies on a thread local" }
]

struct CustomLockType;
struct OtherCustomLockType;
async fn foo() {
  let _x = CustomLockType;
  let _y = OtherCustomLockType;
  baz().await; // Lint violation
}'
$ws.Range("E24").Value = 'if { true } { /* ... */ }

if { let x = somefunc(); x } { /* ... */ }'
$ws.Range("F24").Value = 'if true { /* ... */ }

let res = { let x = somefunc(); x };
if res { /* ... */ }'
$ws.Range("E28").Value = 'if condition {
    1_i64
} else {
    0
};'
$ws.Range("F28").Value = 'i64::from(condition);

or
condition as i64;'
$ws.Range("F29").Value = 'let a: &String = s;'
$ws.Range("E31").Value = 'foo <= i32::MAX as u32;'
$ws.Range("F31").Value = 'i32::try_from(foo).is_ok();'
$ws.Range("E33").Value = 'if x {
    if y {
        // …
    }
}'
$ws.Range("F33").Value = 'if x && y {
    // …
}'
$ws.Range("E35").Value = 'let mut sorted_samples = samples.clone();
sorted_samples.sort();
for sample in &samples { // Oops, meant to use `sorted_samples`.
    println!("{sample}");
}'
$ws.Range("F35").Value = 'let mut sorted_samples = samples.clone();
sorted_samples.sort();
for sample in &sorted_samples {
    println!("{sample}");
}'
$ws.Range("E36").Value = 'fn f(x: u8, y: u8) {
    if x > y {
        a()
    } else if x < y {
        b()
    } else {
        c()
    }
}'
$ws.Range("F36").Value = 'use std::cmp::Ordering;
fn f(x: u8, y: u8) {
     match x.cmp(&y) {
         Ordering::Greater => a(),
         Ordering::Less => b(),
         Ordering::Equal => c()
     }
}'
$ws.Range("E41").Value = 'struct Countdown(u8);

impl Iterator for Countdown {
    // ...
}

let a: Vec<_> = my_iterator.take(1).collect();
let b: Vec<_> = my_iterator.collect();'
$ws.Range("E42").Value = 'macro_rules! print_message {
    () => {
        println!("{}", crate::MESSAGE);
    };
}
pub const MESSAGE: &str = "Hello!";'
$ws.Range("F42").Value = 'macro_rules! print_message {
    () => {
        println!("{}", $crate::MESSAGE);
    };
}
pub const MESSAGE: &str = "Hello!";

Note that if the use of crate is intentional, an allow attribute can be applied to the
macro definition, e.g.:
macro_rules! ok { ... crate::foo ... }'
$ws.Range("E46").Value = 'let mut a: A = Default::default();
a.i = 42;'
$ws.Range("F46").Value = 'let a = A {
    i: 42,
    .. Default::default()
};'
$ws.Range("E47").Value = 'struct S<T> {
    _marker: PhantomData<T>
}

let _: S<i32> = S {
    _marker: PhantomData::default()
};'
$ws.Range("F47").Value = 'struct S<T> {
    _marker: PhantomData<T>
}

let _: S<i32> = S {
    _marker: PhantomData
};'
$ws.Range("F50").Value = 'union Foo {
    a: i32,
    b: u32,
}

fn main() {
    let _x: u32 = unsafe {
        Foo { a: 0_i32 }.b // Now defined behavior, this is just an i32 -> u32 transmute
    };
}'
$ws.Range("F52").Value = 'let x: &i32 = &5;
fun(x);'
$ws.Range("F55").Value = 'struct Foo {
    bar: bool
}'
$ws.Range("E56").Value = 'struct Foo;

impl PartialEq for Foo {
    ...
}'
$ws.Range("E57").Value = 'struct Foo;

impl PartialOrd for Foo {
    ...
}'
$ws.Range("F57").Value = 'struct Foo;

impl PartialOrd for Foo {
    fn partial_cmp(&self, other: &Foo) -> Option<Ordering> {
       Some(self.cmp(other))
    }
}

impl Ord for Foo {
    ...
}

or, if you don''t need a custom ordering:
struct Foo;'
$ws.Range("E58").Value = 'struct Foo;

impl Clone for Foo {
    // ..
}'
$ws.Range("E59").Value = 'use serde::Deserialize;

pub struct Foo {
    // ..
}

impl Foo {
    pub fn new() -> Self {
        // setup here ..
    }

    pub unsafe fn parts() -> (&str, &str) {
        // assumes invariants hold
    }
}'
$ws.Range("E60").Value = 'struct Foo {
    i_am_eq: i32,
    i_am_eq_too: Vec<String>,
}'
$ws.Range("F60").Value = 'struct Foo {
    i_am_eq: i32,
    i_am_eq_too: Vec<String>,
}'
$ws.Range("E61").Value = 'An example clippy.toml configuration:
disallowed-macros = [
    # Can use a string as the path of the disallowed macro.
    "std::print",
    # Can also use an inline table with a `path` key.
    { path = "std::println" },
    # When using an inline table, can add a `reason` for why the macro
    # is disallowed.
    { path = "serde::Serialize", reason = "no serializing" },
]

use serde::Serialize;

// Example code where clippy issues a warning
println!("warns");

// The diagnostic will contain the message "no serializing"
struct Data {
    name: String,
    value: usize,
}'
$ws.Range("E62").Value = 'An example clippy.toml configuration:
disallowed-methods = [
    # Can use a string as the path of the disallowed method.
    "std::boxed::Box::new",
    # Can also use an inline table with a `path` key.
    { path = "std::time::Instant::now" },
    # When using an inline table, can add a `reason` for why the method
    # is disallowed.
    { path = "std::vec::Vec::leak", reason = "no leaking memory" },
]

// Example code where clippy issues a warning
let xs = vec![1, 2, 3, 4];
xs.leak(); // Vec::leak is disallowed in the config.
// The diagnostic contains the message "no leaking memory".

let _now = Instant::now(); // Instant::now is disallowed in the config.

let _box = Box::new(3); // Box::new is disallowed in the config.'
$ws.Range("E65").Value = 'An example clippy.toml configuration:
disallowed-types = [
    # Can use a string as the path of the disallowed type.
    "std::collections::BTreeMap",
    # Can also use an inline table with a `path` key.
    { path = "std::net::TcpListener" },
    # When using an inline table, can add a `reason` for why the type
    # is disallowed.
    { path = "std::net::Ipv4Addr", reason = "no IPv4 allowed" },
]

use std::collections::BTreeMap;
// or its use
let x = std::collections::BTreeMap::new();'
$ws.Range("E73").Value = 'fn simple_double_parens() -> i32 {
    ((0))
}

foo((0));'
$ws.Range("F73").Value = 'fn simple_no_parens() -> i32 {
    0
}

foo(0);'
$ws.Range("E77").Value = '// lib.rs
mod a;
mod b;

// a.rs
mod b;'
$ws.Range("E78").Value = 'if x.is_positive() {
    a();
} else if x.is_negative() {
    b();
}'
$ws.Range("F78").Value = 'if x.is_positive() {
    a();
} else if x.is_negative() {
    b();
} else {
    // We don''t care about zero.
}'
$ws.Range("F80").Value = '
struct Test(!);'
$ws.Range("E82").Value = 'if !map.contains_key(&k) {
    map.insert(k, v);
}'
$ws.Range("F82").Value = 'map.entry(k).or_insert(v);'
$ws.Range("E83").Value = 'enum NonPortable {
    X = 0x1_0000_0000,
    Y = 0,
}'
$ws.Range("F93").Value = 'enum Foo {
    Bar,
    Baz
}'
$ws.Range("F94").Value = 'struct Foo {
    bar: u8,
    baz: String,
}'
$ws.Range("E96").Value = 'writeln!(&mut std::io::stderr(), "foo: {:?}", bar).unwrap();
writeln!(&mut std::io::stdout(), "foo: {:?}", bar).unwrap();'
$ws.Range("F96").Value = 'eprintln!("foo: {:?}", bar);
println!("foo: {:?}", bar);'
$ws.Range("E106").Value = '// &&! looks like a different operator
if foo &&! bar {}'
$ws.Range("F106").Value = 'if foo && !bar {}'
$ws.Range("E109").Value = 'println!("error: {}", format!("something failed at {}", Location::caller()));'
$ws.Range("F109").Value = 'println!("error: something failed at {}", Location::caller());'
$ws.Range("E110").Value = 'println!("error: something failed at {}", Location::caller().to_string());'
$ws.Range("F110").Value = 'println!("error: something failed at {}", Location::caller());'
$ws.Range("E111").Value = 'format!("{}", var);
format!("{v:?}", v = var);
format!("{0} {0}", var);
format!("{0:1$}", var, width);
format!("{:.*}", prec, var);'
$ws.Range("F111").Value = 'format!("{var} {}", 1+2);'
$ws.Range("E117").Value = 'let ptr = Box::into_raw(Box::new(42usize)) as *mut c_void;
let _ = unsafe { Box::from_raw(ptr) };'
$ws.Range("F117").Value = 'let _ = unsafe { Box::from_raw(ptr as *mut usize) };'
$ws.Range("E121").Value = 'if !v.is_empty() {
    a()
} else {
    b()
}'
$ws.Range("F121").Value = 'if v.is_empty() {
    b()
} else {
    a()
}'
$ws.Range("E122").Value = 'let a = if v.is_empty() {
    println!("true!");
    Some(42)
} else {
    None
};'
$ws.Range("F122").Value = 'let a = v.is_empty().then(|| {
    println!("true!");
    42
});'
$ws.Range("E123").Value = 'impl<K: Hash + Eq, V> Serialize for HashMap<K, V> { }

pub fn foo(map: &mut HashMap<i32, i32>) { }

could be rewritten as
impl<K: Hash + Eq, V, S: BuildHasher> Serialize for HashMap<K, V, S> { }

pub fn foo<S: BuildHasher>(map: &mut HashMap<i32, i32, S>) { }'
$ws.Range("E126").Value = 'let mut i: u32 = end - start;

if i != 0 {
    i -= 1;
}'
$ws.Range("F126").Value = 'let mut i: u32 = end - start;

i = i.saturating_sub(1);'
$ws.Range("F127").Value = 'Foo { x, y };'
$ws.Range("F128").Value = '// Index within bounds

x[0];
x[3];'
$ws.Range("F129").Value = '
x.get(2);
x.get(2..100);

y.get(10);
y.get(10..100);'
$ws.Range("E132").Value = 'let infinite_iter = 0..;
[0..].iter().zip(infinite_iter.take_while(|x| *x > 5));'
$ws.Range("E139").Value = 'let time_passed = Instant::now() - Duration::from_secs(5);'
$ws.Range("F139").Value = 'let time_passed = Instant::now().checked_sub(Duration::from_secs(5));'
$ws.Range("E140").Value = 'if x >= y + 1 {}'
$ws.Range("F140").Value = 'if x > y {}'
$ws.Range("E142").Value = 'unsafe {
    std::str::from_utf8_unchecked(b"cl\x82ippy");
}'
$ws.Range("E144").Value = 'mod tests {
    // [...]
}

fn my_function() {
    // [...]
}'
$ws.Range("F144").Value = 'fn my_function() {
    // [...]
}

mod tests {
    // [...]
}'
$ws.Range("F157").Value = 'async fn foo() -> Result<(), ()> {
    Ok(())
}
let _ = foo().await;'
$ws.Range("E162").Value = 'let mut lines = BufReader::new(File::open("some-path")?).lines().filter_map(Result::ok);
// If "some-path" points to a directory, the next statement never terminates:
let first_line: Option<String> = lines.next();'
$ws.Range("F162").Value = 'let mut lines = BufReader::new(File::open("some-path")?).lines().map_while(Result::ok);
let first_line: Option<String> = lines.next();'
$ws.Range("E163").Value = '61864918973511'
$ws.Range("F163").Value = '61_864_918_973_511'
$ws.Range("E165").Value = '618_64_9189_73_511'
$ws.Range("F165").Value = '61_864_918_973_511'
$ws.Range("E169").Value = 'use some_macro;'
$ws.Range("E174").Value = 'if input > max {
    max
} else if input < min {
    min
} else {
    input
}

input.max(min).min(max)

match input {
    x if x > max => max,
    x if x < min => min,
    x => x,
}

let mut x = input;
if x < min { x = min; }
if x > max { x = max; }'
$ws.Range("F174").Value = 'input.clamp(min, max)'
$ws.Range("E176").Value = 'let v = if let Some(v) = w { v } else { return };'
$ws.Range("F176").Value = 'let Some(v) = w else { return };'
$ws.Range("F178").Value = 'struct S {
    pub a: i32,
    pub b: i32,
}

enum E {
    A,
    B,
}

struct T(pub i32, pub i32);'
$ws.Range("E181").Value = 'let newlen = data.len() * std::mem::size_of::<i32>();'
$ws.Range("F181").Value = 'let newlen = std::mem::size_of_val(data);'
$ws.Range("E184").Value = 'let x: Option<String> = do_stuff();
x.map(log_err_msg);
x.map(|msg| log_err_msg(format_msg(msg)));'
$ws.Range("F184").Value = 'let x: Option<String> = do_stuff();
if let Some(msg) = x {
    log_err_msg(msg);
}

if let Some(msg) = x {
    log_err_msg(format_msg(msg));
}'
$ws.Range("E185").Value = 'let x: Result<String, String> = do_stuff();
x.map(log_err_msg);
x.map(|msg| log_err_msg(format_msg(msg)));'
$ws.Range("F185").Value = 'let x: Result<String, String> = do_stuff();
if let Ok(msg) = x {
    log_err_msg(msg);
};
if let Ok(msg) = x {
    log_err_msg(format_msg(msg));
};'
$ws.Range("E187").Value = 'mem::forget(Rc::new(55))'
$ws.Range("E189").Value = 'use std::mem;
///# fn may_panic(v: Vec<i32>) -> Vec<i32> { v }

fn myfunc (v: &mut Vec<i32>) {
    let taken_v = unsafe { mem::replace(v, mem::uninitialized()) };
    let new_v = may_panic(taken_v); // undefined behavior on panic
    mem::forget(mem::replace(v, new_v));
}

The take_mut crate offers a sound solution,
at the cost of either lazily creating a replacement value or aborting
on panic, to ensure that the uninitialized value cannot be observed.'
$ws.Range("E197").Value = 'fn call(service: Service) {
    assert!(service.ready);
}'
$ws.Range("F197").Value = 'fn call(service: Service) {
    assert!(service.ready, "`service.poll_ready()` must be called first to ensure that service is ready to receive requests");
}'
$ws.Range("E198").Value = 'fn new() -> Self {
    Self { random_number: 42 }
}

Could be a const fn:
const fn new() -> Self {
    Self { random_number: 42 }
}'
$ws.Range("E200").Value = 'An example clippy.toml configuration:
enforced-import-renames = [ { path = "serde_json::Value", rename = "JsonValue" }]

use serde_json::Value;'
$ws.Range("E201").Value = 'pub fn foo() {} // missing #[inline]
fn ok() {} // ok

pub trait Bar {
  fn bar(); // ok
  fn def_bar() {} // missing #[inline]
}

struct Baz;
impl Baz {
   fn private() {} // ok
}

impl Bar for Baz {
  fn bar() {} // ok - Baz is not exported
}

pub struct PubBaz;
impl PubBaz {
   fn private() {} // ok
   pub fn not_private() {} // missing #[inline]
}

impl Bar for PubBaz {
   fn bar() {} // missing #[inline]
   fn def_bar() {} // missing #[inline]
}'
$ws.Range("E202").Value = 'trait Trait {
    fn required();

    fn provided() {}
}

impl Trait for Type {
    fn required() { /* ... */ }
}'
$ws.Range("F202").Value = 'trait Trait {
    fn required();

    fn provided() {}
}

impl Trait for Type {
    fn required() { /* ... */ }

    fn provided() { /* ... */ }
}'
$ws.Range("F203").Value = 'let tmp = {
    x = 1;
    1
};
let a = tmp + x;'
$ws.Range("E204").Value = 'let a = b() || panic!() || c();
// `c()` is dead, `panic!()` is only called if `b()` returns `false`
let x = (a, b, c, panic!());
// can simply be replaced by `panic!()`'
$ws.Range("E209").Value = 'debug_assert_eq!(vec![3].pop(), Some(3));

// or

debug_assert!(takes_a_mut_parameter(&mut x));'
$ws.Range("E210").Value = 'let x = Mutex::new(&y);'
$ws.Range("F210").Value = 'let x = AtomicBool::new(y);'
$ws.Range("E211").Value = 'let x = Mutex::new(0usize);'
$ws.Range("F211").Value = 'let x = AtomicUsize::new(0usize);'
$ws.Range("E213").Value = 'let x = &mut &mut y;'
$ws.Range("E214").Value = 'vec.push(&mut value);'
$ws.Range("F214").Value = 'vec.push(&value);'
$ws.Range("E216").Value = 'if x {
    false
} else {
    true
}'
$ws.Range("F216").Value = '!x'
$ws.Range("E218").Value = 'if must_keep(x, y) {
    skip = false;
} else {
    skip = true;
}'
$ws.Range("F218").Value = 'skip = !must_keep(x, y);'
$ws.Range("E220").Value = 'while condition() {
    update_condition();
    if x {
        // ...
    } else {
        continue;
    }
    println!("Hello, world");
}

Could be rewritten as
while condition() {
    update_condition();
    if x {
        // ...
        println!("Hello, world");
    }
}

As another example, the following code
loop {
    if waiting() {
        continue;
    } else {
        // Do something useful
    }
    # break;
}

Could be rewritten as
loop {
    if waiting() {
        continue;
    }
    // Do something useful
    # break;
}'
$ws.Range("E226").Value = 'Point {
    x: 1,
    y: 1,
    z: 1,
    ..zero_point
};'
$ws.Range("F227").Value = 'use std::cmp::Ordering;

let _not_less_or_equal = match a.partial_cmp(&b) {
    None | Some(Ordering::Greater) => true,
    _ => false,
};'
$ws.Range("F231").Value = 'static STATIC_ATOM: AtomicUsize = AtomicUsize::new(15);
STATIC_ATOM.store(9, SeqCst);
assert_eq!(STATIC_ATOM.load(SeqCst), 9); // use a `static` item to refer to the same instance'
$ws.Range("E241").Value = ' fn example(arg_one: u32, arg_two: usize) {}'
$ws.Range("F241").Value = ' extern "C" fn example(arg_one: u32, arg_two: usize) {}'
$ws.Range("E243").Value = 'fn f(a: usize, b: usize) -> usize {
    if a == 0 {
        1
    } else {
        f(a - 1, b + 1)
    }
}'
$ws.Range("F243").Value = 'fn f(a: usize) -> usize {
    if a == 0 {
        1
    } else {
        f(a - 1)
    }
}'
$ws.Range("E245").Value = 'let _ = if let Some(foo) = optional {
    foo
} else {
    5
};
let _ = match optional {
    Some(val) => val + 1,
    None => 5
};
let _ = if let Some(foo) = optional {
    foo
} else {
    let y = do_complicated_function();
    y*y
};

should be
let _ = optional.map_or(5, |foo| foo);
let _ = optional.map_or(5, |val| val + 1);
let _ = optional.map_or_else(||{
    let y = do_complicated_function();
    y*y
}, |foo| foo);'
$ws.Range("E246").Value = 'a + b < a;'
$ws.Range("E256").Value = 'struct TooLarge([u8; 2048]);

fn foo(v: TooLarge) {}'
$ws.Range("F256").Value = 'fn foo(v: &TooLarge) {}'
$ws.Range("E268").Value = 'for i in x..(y+1) {
    // ..
}'
$ws.Range("F268").Value = 'for i in x..=y {
    // ..
}'
$ws.Range("E269").Value = 'for i in x..=(y-1) {
    // ..
}'
$ws.Range("F269").Value = 'for i in x..y {
    // ..
}'
$ws.Range("E275").Value = '{
    let x = Foo::new();
    call(x.clone());
    call(x.clone()); // this can just pass `x`
}

["lorem", "ipsum"].join(" ").to_string();

Path::new("/a/b").join("c").to_path_buf();'
$ws.Range("F290").Value = '// It''s better to have the `#[must_use]` attribute on the method like this:
pub struct Bar;
impl Bar {
    #[must_use]
    pub fn bar(&self) -> Self {
        Self
    }
}

// Or on the type definition like this:
pub struct Bar;
impl Bar {
    pub fn bar(&self) -> Self {
        Self
    }
}'
$ws.Range("E293").Value = 'unsafe { f(x) };'
$ws.Range("F293").Value = 'unsafe { f(x); }'
$ws.Range("E294").Value = 'unsafe { f(x); }'
$ws.Range("F294").Value = 'unsafe { f(x) };'
$ws.Range("E297").Value = 'let x = &x;'
$ws.Range("F297").Value = 'let y = &x; // use different variable name'
$ws.Range("E299").Value = 'let x = y;
let x = z; // shadows the earlier binding'
$ws.Range("F299").Value = 'let x = y;
let w = z; // use different variable name'
$ws.Range("E303").Value = 'const SIZE: usize = 128;
let x = [2u8; SIZE];
let mut y = [2u8; SIZE];
unsafe { copy_nonoverlapping(x.as_ptr(), y.as_mut_ptr(), size_of::<u8>() * SIZE) };'
$ws.Range("E305").Value = 'let mut vec1 = Vec::with_capacity(len);
vec1.resize(len, 0);

let mut vec1 = Vec::with_capacity(len);
vec1.resize(vec1.capacity(), 0);

let mut vec2 = Vec::with_capacity(len);
vec2.extend(repeat(0).take(len));'
$ws.Range("F305").Value = 'let mut vec1 = vec![0; len];
let mut vec2 = vec![0; len];'
$ws.Range("F307").Value = 'use alloc::vec::Vec;'
$ws.Range("E308").Value = 'use alloc::slice::from_ref;'
$ws.Range("F319").Value = '// same as above except:
impl PartialEq for Vec3 {
    fn eq(&self, other: &Self) -> bool {
        // Note we now compare other.x to self.x
        self.x == other.x && self.y == other.y && self.z == other.z
    }
}'
$ws.Range("E324").Value = 'a = b;
b = a;'
$ws.Range("F324").Value = 'std::mem::swap(&mut a, &mut b);'
$ws.Range("E328").Value = 'fn my_cool_test() {
    // [...]
}

mod tests {
    // [...]
}'
$ws.Range("F328").Value = 'mod tests {
    #[test]
    fn my_cool_test() {
        // [...]
    }
}'
$ws.Range("E329").Value = 'let is_digit = c.to_digit(radix).is_some();'
$ws.Range("F329").Value = 'let is_digit = c.is_digit(radix);'
$ws.Range("F330").Value = 'struct MoreOftenUseful {
    some_field: usize,
    last: [u32; 0],
}'
$ws.Range("E355").Value = 'if option.is_some() {
    do_something_with(option.unwrap())
}'
$ws.Range("F355").Value = 'if let Some(value) = option {
    do_something_with(value)
}'
$ws.Range("E356").Value = 'if option.is_none() {
    do_something_with(option.unwrap())
}

This code will always panic. The if condition should probably be inverted.'
$ws.Range("F361").Value = 'foo(&[1, 2]);'
$ws.Range("E363").Value = 'use std::cmp::Ordering::*;

foo(Less);'
$ws.Range("F363").Value = 'use std::cmp::Ordering;

foo(Ordering::Less)'
$ws.Range("E366").Value = 'print!("Hello {}!\n", name);

use println!() instead
println!("Hello {}!", name);'
$ws.Range("E369").Value = 'println!("{:?}", foo);'
$ws.Range("E371").Value = 'writeln!(buf, "");'
$ws.Range("F371").Value = 'writeln!(buf);'
$ws.Range("E372").Value = 'write!(buf, "Hello {}!\n", name);'
$ws.Range("F372").Value = 'writeln!(buf, "Hello {}!", name);'
$ws.Range("E373").Value = 'writeln!(buf, "{}", "foo");'
$ws.Range("F373").Value = 'writeln!(buf, "foo");'
$ws.Range("E375").Value = 'fn unique_words(text: &str) -> HashMap<&str, ()> {
    todo!();
}'
$ws.Range("F375").Value = 'fn unique_words(text: &str) -> HashSet<&str> {
    todo!();
}'
$ws.Range("E376").Value = '[package]
name = "clippy"
version = "0.0.212"
repository = "https://github.com/rust-lang/rust-clippy"
readme = "README.md"
license = "MIT OR Apache-2.0"
keywords = ["clippy", "lint", "plugin"]
categories = ["development-tools", "development-tools::cargo-plugins"]

Should include a description field like:
[package]
name = "clippy"
version = "0.0.212"
description = "A bunch of helpful lints to avoid common pitfalls in Rust"
repository = "https://github.com/rust-lang/rust-clippy"
readme = "README.md"
license = "MIT OR Apache-2.0"
keywords = ["clippy", "lint", "plugin"]
categories = ["development-tools", "development-tools::cargo-plugins"]'
$ws.Range("E377").Value = '[features]
default = ["use-abc", "with-def", "ghi-support"]
use-abc = []  // redundant
with-def = []   // redundant
ghi-support = []   // redundant'
$ws.Range("E378").Value = '[features]
default = []
no-abc = []
not-def = []'
$ws.Range("E379").Value = '[dependencies]
ctrlc = "=3.1.0"
ansi_term = "=0.11.0"'
$ws.Range("F383").Value = 'fn as_u8(x: u64) -> u8 {
    if let Ok(x) = u8::try_from(x) {
        x
    } else {
        todo!();
    }
}
// Or
fn as_u16(x: u64) -> u16 {
    x as u16
}'
$ws.Range("F388").Value = 'let _ = fun as usize;'
$ws.Range("E403").Value = 'fn foo(x: u32, y: u32, name: &str, c: Color, w: f32, h: f32, a: f32, b: f32) {
    // ..
}'
$ws.Range("E406").Value = 'fn useless() { }'
$ws.Range("E407").Value = 'fn double_must_use() -> Result<(), ()> {
    unimplemented!();
}'
$ws.Range("E409").Value = 'pub fn read_u8() -> Result<u8, ()> { Err(()) }

should become
use std::fmt;

pub struct EndOfStream;

impl fmt::Display for EndOfStream {
    fn fmt(&self, f: &mut fmt::Formatter<''_>) -> fmt::Result {
        write!(f, "End of Stream")
    }
}

impl std::error::Error for EndOfStream { }

pub fn read_u8() -> Result<u8, EndOfStream> { Err(EndOfStream) }
///# fn main() {
///#     read_u8().unwrap();
///# }

Note that there are crates that simplify creating the error type, e.g.
thiserror.'
$ws.Range("E413").Value = 'for i in 0..src.len() {
    dst[i + 64] = src[i];
}'
$ws.Range("F413").Value = 'dst[64..(src.len() + 64)].clone_from_slice(&src[..]);'
$ws.Range("E415").Value = '// with `y` a `Vec` or slice:
for x in y.iter() {
    // ..
}'
$ws.Range("F415").Value = 'for x in &y {
    // ..
}'
$ws.Range("E416").Value = '// with `y` a `Vec` or slice:
for x in y.into_iter() {
    // ..
}

can be rewritten to
for x in y {
    // ..
}'
$ws.Range("E418").Value = 'loop {
    let x = match y {
        Some(x) => x,
        None => break,
    };
    // .. do something with x
}
// is easier written as
while let Some(x) = y {
    // .. do something with x
};'
$ws.Range("E419").Value = 'let mut i = 0;
for item in &v {
    bar(i, *item);
    i += 1;
}'
$ws.Range("F419").Value = 'for (i, item) in v.iter().enumerate() { bar(i, *item); }'
$ws.Range("E432").Value = 'match x {
    Some(ref foo) => bar(foo),
    _ => (),
}'
$ws.Range("F432").Value = 'if let Some(ref foo) = x {
    bar(foo);
}'
$ws.Range("E433").Value = 'Using match:
match x {
    Some(ref foo) => bar(foo),
    _ => bar(&other_ref),
}

Using if let with else:
if let Some(ref foo) = x {
    bar(foo);
} else {
    bar(&other_ref);
}'
$ws.Range("E435").Value = 'let condition: bool = true;
match condition {
    true => foo(),
    false => bar(),
}'
$ws.Range("F435").Value = 'let condition: bool = true;
if condition {
    foo();
} else {
    bar();
}'
$ws.Range("E439").Value = 'match x {
    Foo::A(_) => {},
    _ => {},
}'
$ws.Range("F439").Value = 'match x {
    Foo::A(_) => {},
    Foo::B(_) => {},
}'
$ws.Range("E440").Value = 'match x {
    Foo::A => {},
    Foo::B => {},
    _ => {},
}'
$ws.Range("F440").Value = 'match x {
    Foo::A => {},
    Foo::B => {},
    Foo::C => {},
}'
$ws.Range("E441").Value = 'match s {
    "a" => {},
    "bar" | _ => {},
}'
$ws.Range("F441").Value = 'match s {
    "a" => {},
    _ => {},
}'
$ws.Range("E443").Value = 'match (a, b) {
    (c, d) => {
        // useless match
    }
}'
$ws.Range("F443").Value = 'let (c, d) = (a, b);'
$ws.Range("E444").Value = 'let a = A { a: 5 };

match a {
    A { a: 5, .. } => {},
    _ => {},
}'
$ws.Range("F444").Value = 'match a {
    A { a: 5 } => {},
    _ => {},
}'
$ws.Range("E445").Value = 'if let Ok(_) = Ok::<i32, i32>(42) {}
if let Err(_) = Err::<i32, i32>(42) {}
if let None = None::<()> {}
if let Some(_) = Some(42) {}
if let Poll::Pending = Poll::Pending::<()> {}
if let Poll::Ready(_) = Poll::Ready(42) {}
if let IpAddr::V4(_) = IpAddr::V4(Ipv4Addr::LOCALHOST) {}
if let IpAddr::V6(_) = IpAddr::V6(Ipv6Addr::LOCALHOST) {}
match Ok::<i32, i32>(42) {
    Ok(_) => true,
    Err(_) => false,
};'
$ws.Range("F445").Value = 'if Ok::<i32, i32>(42).is_ok() {}
if Err::<i32, i32>(42).is_err() {}
if None::<()>.is_none() {}
if Some(42).is_some() {}
if Poll::Pending::<()>.is_pending() {}
if Poll::Ready(42).is_ready() {}
if IpAddr::V4(Ipv4Addr::LOCALHOST).is_ipv4() {}
if IpAddr::V6(Ipv6Addr::LOCALHOST).is_ipv6() {}
Ok::<i32, i32>(42).is_ok();'
$ws.Range("E452").Value = 'match &*text.to_ascii_lowercase() {
    "foo" => {},
    "Bar" => {},
    _ => {},
}'
$ws.Range("F452").Value = 'match &*text.to_ascii_lowercase() {
    "foo" => {},
    "bar" => {},
    _ => {},
}'
$ws.Range("E453").Value = 'let mutex = Mutex::new(State {});

match mutex.lock().unwrap().foo() {
    true => {
        mutex.lock().unwrap().bar(); // Deadlock!
    }
    false => {}
};

println!("All done!");'
$ws.Range("F453").Value = 'let mutex = Mutex::new(State {});

let is_foo = mutex.lock().unwrap().foo();
match is_foo {
    true => {
        mutex.lock().unwrap().bar();
    }
    false => {}
};

println!("All done!");'
$ws.Range("E459").Value = 'vec.iter().cloned().take(10);
vec.iter().cloned().last();'
$ws.Range("F459").Value = 'vec.iter().take(10).cloned();
vec.iter().last().cloned();'
$ws.Range("E461").Value = 'option.unwrap();
result.unwrap();'
$ws.Range("F461").Value = 'option?;

// or

result?;'
$ws.Range("E462").Value = 'option.expect("one");
result.expect("one");'
$ws.Range("F462").Value = 'option?;

// or

result?;'
$ws.Range("E463").Value = 'struct X;
impl X {
    fn add(&self, other: &X) -> X {
        // ..
    }
}'
$ws.Range("E464").Value = 'impl X {
    fn as_str(self) -> &''static str {
        // ..
    }
}'
$ws.Range("E465").Value = 'x.ok().expect("why did I do this again?");'
$ws.Range("F465").Value = 'x.expect("why did I do this again?");'
$ws.Range("E467").Value = 'x.unwrap_or_else(Default::default);
x.unwrap_or_else(u32::default);'
$ws.Range("F467").Value = 'x.unwrap_or_default();'
$ws.Range("E468").Value = 'option.map(|a| a + 1).unwrap_or(0);
result.map(|a| a + 1).unwrap_or_else(some_function);'
$ws.Range("F468").Value = 'option.map_or(0, |a| a + 1);
result.map_or_else(some_function, |a| a + 1);'
$ws.Range("E469").Value = 'opt.map_or(None, |a| Some(a + 1));'
$ws.Range("F469").Value = 'opt.and_then(|a| Some(a + 1));'
$ws.Range("E470").Value = 'assert_eq!(Some(1), r.map_or(None, Some));'
$ws.Range("F470").Value = 'assert_eq!(Some(1), r.ok());'
$ws.Range("E471").Value = 'let _ = opt().and_then(|s| Some(s.len()));
let _ = res().and_then(|s| if s.len() == 42 { Ok(10) } else { Ok(20) });
let _ = res().or_else(|s| if s.len() == 42 { Err(10) } else { Err(20) });'
$ws.Range("F471").Value = 'let _ = opt().map(|s| s.len());
let _ = res().map(|s| if s.len() == 42 { 10 } else { 20 });
let _ = res().map_err(|s| if s.len() == 42 { 10 } else { 20 });'
$ws.Range("E472").Value = 'vec.iter().filter(|x| **x == 0).next();'
$ws.Range("F472").Value = 'vec.iter().find(|x| **x == 0);'
$ws.Range("E473").Value = 'vec.iter().skip_while(|x| **x == 0).next();'
$ws.Range("F473").Value = 'vec.iter().find(|x| **x != 0);'
$ws.Range("F474").Value = 'vec.iter().flat_map(|x| x.iter());
opt.and_then(|x| Some(x * 2));'
$ws.Range("E475").Value = '(0_i32..10)
    .filter(|n| n.checked_add(1).is_some())
    .map(|n| n.checked_add(1).unwrap());'
$ws.Range("F475").Value = '(0_i32..10).filter_map(|n| n.checked_add(1));'
$ws.Range("E478").Value = 'iter.flat_map(|x| x);

Can be written as
iter.flatten();'
$ws.Range("E479").Value = 'let vec = vec![1];
vec.iter().find(|x| **x == 0).is_some();

"hello world".find("world").is_none();'
$ws.Range("F479").Value = 'let vec = vec![1];
vec.iter().any(|x| *x == 0);

!"hello world".contains("world");'
$ws.Range("E481").Value = 'foo.unwrap_or(String::from("empty"));'
$ws.Range("F481").Value = 'foo.unwrap_or_else(|| String::from("empty"));'
$ws.Range("E482").Value = '// Result
let value = result.or::<Error>(Ok(fallback)).unwrap();

// Option
let value = option.or(Some(fallback)).unwrap();'
$ws.Range("F482").Value = '// Result
let value = result.unwrap_or(fallback);

// Option
let value = option.unwrap_or(fallback);'
$ws.Range("E483").Value = 'foo.expect(&format!("Err {}: {}", err_code, err_msg));

// or

foo.expect(format!("Err {}: {}", err_code, err_msg).as_str());'
$ws.Range("F483").Value = 'foo.unwrap_or_else(|| panic!("Err {}: {}", err_code, err_msg));'
$ws.Range("E485").Value = 'let x = Rc::new(1);

x.clone();'
$ws.Range("F485").Value = 'Rc::clone(&x);'
$ws.Range("E487").Value = 'In an impl block:
impl Foo {
    fn new() -> NotAFoo {
    }
}

struct Bar(Foo);
impl Foo {
    // Bad. The type name must contain `Self`
    fn new() -> Bar {
    }
}

impl Foo {
    // Good. Return type contains `Self`
    fn new() -> Result<Foo, FooError> {
    }
}

Or in a trait definition:
pub trait Trait {
    // Bad. The type name must contain `Self`
    fn new();
}

pub trait Trait {
    // Good. Return type contains `Self`
    fn new() -> Self;
}'
$ws.Range("E491").Value = 'let x = s.iter().nth(0);'
$ws.Range("F491").Value = 'let x = s.iter().next();'
$ws.Range("E494").Value = 'let mut foo = vec![0, 1, 2, 3];
let bar: HashSet<usize> = foo.drain(..).collect();'
$ws.Range("F494").Value = 'let foo = vec![0, 1, 2, 3];
let bar: HashSet<usize> = foo.into_iter().collect();'
$ws.Range("E500").Value = 'name.chars().last() == Some(''_'') || name.chars().next_back() == Some(''-'');'
$ws.Range("F500").Value = 'name.ends_with(''_'') || name.ends_with(''-'');'
$ws.Range("E501").Value = 'let x: &[i32] = &[1, 2, 3, 4, 5];
do_stuff(x.as_ref());'
$ws.Range("F501").Value = 'let x: &[i32] = &[1, 2, 3, 4, 5];
do_stuff(x);'
$ws.Range("E502").Value = '(0..3).fold(false, |acc, x| acc || x > 2);'
$ws.Range("E505").Value = '(&vec).into_iter();'
$ws.Range("F505").Value = '(&vec).iter();'
$ws.Range("E508").Value = 'let add = x.checked_add(y).unwrap_or(u32::MAX);
let sub = x.checked_sub(y).unwrap_or(u32::MIN);'
$ws.Range("F508").Value = 'let add = x.saturating_add(y);
let sub = x.saturating_sub(y);'
$ws.Range("E510").Value = 'let metadata = std::fs::metadata("foo.txt")?;
let filetype = metadata.file_type();

if filetype.is_file() {
    // read file
}'
$ws.Range("F510").Value = 'let metadata = std::fs::metadata("foo.txt")?;
let filetype = metadata.file_type();

if !filetype.is_dir() {
    // read file
}'
$ws.Range("E511").Value = 'opt.as_ref().map(String::as_str)

Can be written as
opt.as_deref()'
$ws.Range("E512").Value = 'a[2..].iter().next();
b.iter().next();'
$ws.Range("F512").Value = 'a.get(2);
b.get(0);'
$ws.Range("E513").Value = 'string.insert_str(0, "R");
string.push_str("R");'
$ws.Range("F513").Value = 'string.insert(0, ''R'');
string.push(''R'');'
$ws.Range("E518").Value = 'iter.filter_map(|x| x);'
$ws.Range("F518").Value = 'iter.flatten();'
$ws.Range("E520").Value = '"Hello".bytes().nth(3);'
$ws.Range("F520").Value = '"Hello".as_bytes().get(3);'
$ws.Range("E522").Value = 'let some_vec = vec![0, 1, 2, 3];

some_vec.iter().count();
&some_vec[..].iter().count();'
$ws.Range("E523").Value = 'let s = "Hello world!";
let cow = Cow::Borrowed(s);

let data = cow.to_owned();
assert!(matches!(data, Cow::Borrowed(_)))'
$ws.Range("F523").Value = 'let s = "Hello world!";
let cow = Cow::Borrowed(s);

let data = cow.clone();
assert!(matches!(data, Cow::Borrowed(_)))

or
let s = "Hello world!";
let cow = Cow::Borrowed(s);

let _data: String = cow.into_owned();'
$ws.Range("E524").Value = 'for x in s.splitn(1, ":") {
    // ..
}'
$ws.Range("F524").Value = 'for x in s.splitn(2, ":") {
    // ..
}'
$ws.Range("E537").Value = 'let count = vec.iter().filter(|x| **x == 0u8).count();'
$ws.Range("F537").Value = 'let count = bytecount::count(&vec, 0u8);'
$ws.Range("E543").Value = 'Before:
use std::fmt;

enum Error {
    Indivisible,
    Remainder(u8),
}

impl fmt::Display for Error {
    fn fmt(&self, f: &mut fmt::Formatter<''_>) -> fmt::Result {
        match self {
            Error::Indivisible => write!(f, "could not divide input by three"),
            Error::Remainder(remainder) => write!(
                f,
                "input is not divisible by three, remainder = {}",
                remainder
            ),
        }
    }
}

impl std::error::Error for Error {}

fn divisible_by_3(input: &str) -> Result<(), Error> {
    input
        .parse::<i32>()
        .map_err(|_| Error::Indivisible)
        .map(|v| v % 3)
        .and_then(|remainder| {
            if remainder == 0 {
                Ok(())
            } else {
                Err(Error::Remainder(remainder as u8))
            }
        })
}'
$ws.Range("F543").Value = 'use std::{fmt, num::ParseIntError};

enum Error {
   Indivisible(ParseIntError),
   Remainder(u8),
}

impl fmt::Display for Error {
   fn fmt(&self, f: &mut fmt::Formatter<''_>) -> fmt::Result {
       match self {
           Error::Indivisible(_) => write!(f, "could not divide input by three"),
           Error::Remainder(remainder) => write!(
               f,
               "input is not divisible by three, remainder = {}",
               remainder
           ),
       }
   }
}

impl std::error::Error for Error {
   fn source(&self) -> Option<&(dyn std::error::Error + ''static)> {
       match self {
           Error::Indivisible(source) => Some(source),
           _ => None,
       }
   }
}

fn divisible_by_3(input: &str) -> Result<(), Error> {
   input
       .parse::<i32>()
       .map_err(Error::Indivisible)
       .map(|v| v % 3)
       .and_then(|remainder| {
           if remainder == 0 {
               Ok(())
           } else {
               Err(Error::Remainder(remainder as u8))
           }
       })
}'
$ws.Range("E547").Value = 'let _ = x.iter().zip(0..x.len());'
$ws.Range("F547").Value = 'let _ = x.iter().enumerate();'
$ws.Range("E550").Value = 'match my_enum {
	Empty => ().hash(&mut state),
	WithValue(x) => x.hash(&mut state),
}'
$ws.Range("F550").Value = 'match my_enum {
	Empty => 0_u8.hash(&mut state),
	WithValue(x) => x.hash(&mut state),
}'
$ws.Range("E551").Value = 'vec.sort_by(|a, b| a.foo().cmp(&b.foo()));'
$ws.Range("F551").Value = 'vec.sort_by_key(|a| a.foo());'
$ws.Range("E553").Value = 'let mut f = File::open("foo.txt").unwrap();
let mut bytes = Vec::new();
f.read_to_end(&mut bytes).unwrap();

Can be written more concisely as
let mut bytes = fs::read("foo.txt").unwrap();'
$ws.Range("E554").Value = 'let map: HashMap<u32, u32> = HashMap::new();
let values = map.iter().map(|(_, value)| value).collect::<Vec<_>>();'
$ws.Range("F554").Value = 'let map: HashMap<u32, u32> = HashMap::new();
let values = map.values().collect::<Vec<_>>();'
$ws.Range("E556").Value = 'fn foo<T: io::Seek>(t: &mut T) {
    t.seek(io::SeekFrom::Start(0));
}'
$ws.Range("F556").Value = 'fn foo<T: io::Seek>(t: &mut T) {
    t.rewind();
}'
$ws.Range("E557").Value = 'let len = iterator.collect::<Vec<_>>().len();'
$ws.Range("F557").Value = 'let len = iterator.count();'
$ws.Range("E560").Value = 'foo.iter().rev().next();'
$ws.Range("F560").Value = 'foo.iter().next_back();'
$ws.Range("E561").Value = 'let f = Foo { a: 0, b: 0, c: 0 };

match f {
    Foo { a: _, b: 0, .. } => {},
    Foo { a: _, b: _, c: _ } => {},
}'
$ws.Range("F561").Value = 'let f = Foo { a: 0, b: 0, c: 0 };

match f {
    Foo { b: 0, .. } => {},
    Foo { .. } => {},
}'
$ws.Range("E564").Value = '0x1a9BAcD'
$ws.Range("F564").Value = '0x1A9BACD'
$ws.Range("E565").Value = '123832i32'
$ws.Range("F565").Value = '123832_i32'
$ws.Range("E566").Value = '123832_i32'
$ws.Range("F566").Value = '123832i32'
$ws.Range("E567").Value = 'In Rust:
fn main() {
    let a = 0123;
    println!("{}", a);
}

prints 123, while in C:

int main() {
    int a = 0123;
    printf("%d\n", a);
}

prints 83 (as 83 == 0o123 while 123 == 0o173).'
$ws.Range("E569").Value = 'match v {
    Some(x) => (),
    y @ _ => (),
}'
$ws.Range("F569").Value = 'match v {
    Some(x) => (),
    y => (),
}'
$ws.Range("E570").Value = 'match t {
    TupleStruct(0, .., _) => (),
    _ => (),
}'
$ws.Range("F570").Value = 'match t {
    TupleStruct(0, ..) => (),
    _ => (),
}'
$ws.Range("E573").Value = 'a + 1.0;'
$ws.Range("E576").Value = 'if (x & 1 == 2) { }'
$ws.Range("E577").Value = 'if (x | 1 > 3) {  }'
$ws.Range("E578").Value = 'if x & 0b1111 == 0 { }'
$ws.Range("E579").Value = 'if x == y || x < y {}'
$ws.Range("F579").Value = 'if x <= y {}'
$ws.Range("E580").Value = 'let micros = duration.subsec_nanos() / 1_000;
let millis = duration.subsec_nanos() / 1_000_000;'
$ws.Range("F580").Value = 'let micros = duration.subsec_micros();
let millis = duration.subsec_millis();'
$ws.Range("E581").Value = 'if x + 1 == x + 1 {}

// or

assert_eq!(a, a);'
$ws.Range("E585").Value = 'x / 1 + 0 * 1 - 0 | 0;'
$ws.Range("E587").Value = 'if x == f32::NAN { }'
$ws.Range("F587").Value = 'if x.is_nan() { }'
$ws.Range("E588").Value = 'if x.to_owned() == y {}'
$ws.Range("F588").Value = 'if x == y {}'
$ws.Range("F589").Value = 'let error_margin = f64::EPSILON; // Use an epsilon for comparison
// Or, if Rust <= 1.42, use `std::f64::EPSILON` constant instead.
// let error_margin = std::f64::EPSILON;
if (y - 1.23f64).abs() < error_margin { }
if (y - x).abs() > error_margin { }'
$ws.Range("F590").Value = 'let error_margin = f64::EPSILON; // Use an epsilon for comparison
// Or, if Rust <= 1.42, use `std::f64::EPSILON` constant instead.
// let error_margin = std::f64::EPSILON;
if (x - ONE).abs() < error_margin { }'
$ws.Range("E591").Value = 'let a = x % 1;
let a = x % -1;'
$ws.Range("E598").Value = 'unsafe { std::mem::transmute::<*const [i32], *const [u16]>(p) };'
$ws.Range("F598").Value = 'p as *const [u16];'
$ws.Range("E605").Value = 'let _non_zero: NonZeroU32 = unsafe { std::mem::transmute(123) };'
$ws.Range("F605").Value = 'let _non_zero = unsafe { NonZeroU32::new_unchecked(123) };'
$ws.Range("F610").Value = 'struct Foo<T>(u32, T);
let _ = unsafe { core::mem::transmute::<Foo<u32>, Foo<i32>>(Foo(0u32, 0u32)) };'
$ws.Range("E616").Value = 'let x: LinkedList<usize> = LinkedList::new();'
$ws.Range("E618").Value = 'fn foo(bar: Rc<&usize>) {}'
$ws.Range("E619").Value = 'fn foo(interned: Rc<String>) { ... }'
$ws.Range("E620").Value = 'struct Foo {
    inner: Rc<Vec<Vec<Box<(u32, u32, u32, u32)>>>>,
}'
$ws.Range("E623").Value = 'if {
    foo();
} == {
    bar();
} {
    baz();
}

is equal to
{
    foo();
    bar();
    baz();
}

For asserts:
assert_eq!({ foo(); }, { bar(); });

will always succeed'
$ws.Range("E625").Value = '[1, 2, 3].into_iter().for_each(|n| { *n; });'
$ws.Range("E627").Value = 'struct Foo {
    x: Box<isize>,
}'
$ws.Range("E629").Value = 'fn main() {
    unsafe {

    }
}'
$ws.Range("E630").Value = 'pub fn foo() {}'
$ws.Range("E631").Value = 'pub struct Foo {
    pub field: i32
}'
$ws.Range("E632").Value = 'pub struct Foo;'
$ws.Range("E633").Value = '// edition 2015
pub trait Foo {
    fn foo(usize);
}
fn main() {}'
$ws.Range("E634").Value = 'const FOO: i32 = 5;'
$ws.Range("E635").Value = 'fn foo<T>(t: T) {

}'
$ws.Range("E638").Value = 'async fn foo() {}'
$ws.Range("E639").Value = 'mod foo {
    pub mod bar {

    }
}'
$ws.Range("E641").Value = 'pub struct A where i32: Copy;'
$ws.Range("E644").Value = '// edition 2015
fn dyn() {}'
$ws.Range("E645").Value = ''
$ws.Range("E646").Value = 'unsafe {
    let x: &''static i32 = std::mem::zeroed();
}'
$ws.Range("E648").Value = 'use std::ptr;
unsafe {
    let x = &*ptr::null::<i32>();
    let x = ptr::addr_of!(*ptr::null::<i32>());
    let x = *(0 as *const i32);
}'
$ws.Range("E649").Value = 'use std::arch::asm;

fn main() {
    unsafe {
        asm!("foo: bar");
    }
}'
$ws.Range("E651").Value = '
use core::ops::Deref;

trait A {}
trait B: A {}
impl<''a> Deref for dyn ''a + B {
    type Target = dyn A;
    fn deref(&self) -> &Self::Target {
        todo!()
    }
}

fn take_a(_: &dyn A) { }

fn take_b(b: &dyn B) {
    take_a(b);
}'
$ws.Range("E652").Value = 'let mut lock_guard = mutex.lock();
std::mem::drop(&lock_guard); // Should have been drop(lock_guard), mutex
// still locked
operation_that_requires_mutex_to_be_unlocked();'
$ws.Range("E656").Value = 'core::mem::discriminant::<i32>(&123);'
$ws.Range("E658").Value = 'fn main() {
    println!("{:?}", ''‮'');
}'
$ws.Range("E662").Value = 'let c_str = CString::new("foo").unwrap().as_ptr();'
$ws.Range("E663").Value = 'trait A {}
trait B {}

trait C: A + B {}'
$ws.Range("E667").Value = 'fn main() {
    let föö = 1;
}'
$ws.Range("E668").Value = 'const µ: f64 = 0.000001;'
$ws.Range("E672").Value = 'struct Foo;
let foo = &Foo;
let clone: &Foo = foo.clone();'
$ws.Range("E673").Value = 'struct Foo;
let foo = &&Foo;
let clone: &Foo = foo.clone();'
$ws.Range("E674").Value = '
trait Duh {}

impl Duh for i32 {}

trait Trait {
    type Assoc: Duh;
}

struct Struct;

impl<F: Duh> Trait for F {
    type Assoc = F;
}

type Tait = impl Sized;

fn test() -> impl Trait<Assoc = Tait> {
    42
}

In this example, test declares that the associated type Assoc for
impl Trait is impl Sized, which does not satisfy the Send bound
on the associated type.
Although the hidden type, i32 does satisfy this bound, we do not
consider the return type to be well-formed with this lint. It can be
fixed by changing Tait = impl Sized into Tait = impl Sized + Send.'
$ws.Range("E680").Value = 'enum En {
    V0(u8),
    VBig([u8; 1024]),
}'
$ws.Range("E682").Value = 'pub extern "C" fn str_type(p: &str) { }'
$ws.Range("E683").Value = 'let atom = AtomicU8::new(0);
let value = atom.load(Ordering::Release);'
$ws.Range("E685").Value = 'fn foo<T>() -> T { panic!() }

fn main() {
    foo::<usize>();
}'
$ws.Range("E689").Value = 'use test::{A};

pub mod test {
    pub struct A;
}'
$ws.Range("E691").Value = '
fn main() {}'
$ws.Range("E692").Value = 'fn foo() {}'
$ws.Range("E693").Value = 'enum Foo {
    Variant1,
}'
$ws.Range("E694").Value = '
macro_rules! foo {
    () => {};
    ($( $i:ident = $($j:ident),+ );*) => { $( $( $i = $k; )+ )* };
}

fn main() {
    foo!();
}'
$ws.Range("E697").Value = 'let x = 1 / 0;'
$ws.Range("E699").Value = '
struct SyncThing {}

async fn yield_now() {}

pub async fn uhoh() {
    let guard = SyncThing {};
    yield_now().await;
}'
$ws.Range("E700").Value = 'extern crate proc_macro;'
$ws.Range("E701").Value = '
This will produce:
error: external crate `regex` unused in `lint_example`: remove the dependency or add `use regex as _;`
  |
note: the lint level is defined here
 --> src/lib.rs:1:9
  |
1 | #![deny(unused_crate_dependencies)]
  |         ^^^^^^^^^^^^^^^^^^^^^^^^^'
$ws.Range("E702").Value = 'mod foo {
    pub fn bar() {}
}

fn main() {
    use foo::bar;
    foo::bar();
}'
$ws.Range("E703").Value = ''
$ws.Range("E704").Value = '
let x = 10;
println!("{}", x);'
$ws.Range("E708").Value = ''
$ws.Range("E709").Value = 'struct S(i32, i32, i32);
let s = S(1, 2, 3);
let _ = (s.0, s.2);'
$ws.Range("E715").Value = 'macro_rules! unused_empty {
    (hello) => { println!("Hello, world!") }; // This rule is unused
    () => { println!("empty") }; // This rule is used
}

fn main() {
    unused_empty!(hello);
}'
$ws.Range("E716").Value = 'fn foo() {}'
$ws.Range("E718").Value = 'fn main() {}'
$ws.Range("E719").Value = 'fn main() {}'
$ws.Range("E720").Value = 'let x: &u32 = &42;
let y = x as *const u32;'
$ws.Range("E721").Value = 'let x = 42_i32 as i32;'
$ws.Range("E722").Value = 'struct SemiPriv;

mod m1 {
    struct Priv;
    impl super::SemiPriv {
        pub fn f(_: Priv) {}
    }
}'
$ws.Range("E723").Value = 'const FOO: () = unsafe {
    let x = &[0_u8; 4];
    let y = x.as_ptr().cast::<u32>();
    let mut z = 123;
    y.copy_to_nonoverlapping(&mut z, 1); // the address of a `u8` array is unknown
    // and thus we don''t know if it is aligned enough for copying a `u32`.
};'
$ws.Range("E727").Value = ''
$ws.Range("E734").Value = 'fn foo() {}

fn bar() {
    foo();
}'
$ws.Range("E738").Value = '
fn foo<''a>(x: &''a u32) {}'
$ws.Range("E739").Value = '
pub fn foo<''a>() {}'
$ws.Range("E741").Value = 'struct Foo<''a> {
    x: &''a u32
}

fn foo(x: &Foo) {
}'
$ws.Range("E743").Value = '
mod foo {
    pub fn bar() {}
}

fn main() {
    ::foo::bar();
}'
$ws.Range("E749").Value = '// foo.rs

extern crate proc_macro;

use proc_macro::*;

pub fn foo1(a: TokenStream) -> TokenStream {
    drop(a);
    "mod __bar { static mut BAR: Option<Something> = None; }".parse().unwrap()
}

// bar.rs
extern crate foo;

struct Something;

struct Another;

fn main() {}

This will produce:
warning: cannot find type `Something` in this scope
 --> src/main.rs:8:10
  |
8 | #[derive(Foo)]
  |          ^^^ names from parent modules are not accessible without an explicit import
  |
  = note: `#[warn(proc_macro_derive_resolution_fallback)]` on by default
  = warning: this was previously accepted by the compiler but is being phased out; it will become a hard error in a future release!
  = note: for more information, see issue #50504 <https://github.com/rust-lang/rust/issues/50504>'
$ws.Range("E750").Value = '
extern crate serde_json;

fn main() {
    let _ = json!{{}};
}

This will produce:
error: deprecated `#[macro_use]` attribute used to import macros should be replaced at use sites with a `use` item to import the macro instead
 --> src/main.rs:3:1
  |
3 | #[macro_use]
  | ^^^^^^^^^^^^
  |
note: the lint level is defined here
 --> src/main.rs:1:9
  |
1 | #![deny(macro_use_extern_crate)]
  |         ^^^^^^^^^^^^^^^^^^^^^^'
$ws.Range("E752").Value = '
struct SharedRef<''a, T>
where
    T: ''a,
{
    data: &''a T,
}'
$ws.Range("E753").Value = '
struct NoDerive(i32);
impl PartialEq for NoDerive { fn eq(&self, _: &Self) -> bool { false } }
impl Eq for NoDerive { }
struct WrapParam<T>(T);
const WRAP_INDIRECT_PARAM: & &WrapParam<NoDerive> = & &WrapParam(NoDerive(0));
fn main() {
    match WRAP_INDIRECT_PARAM {
        WRAP_INDIRECT_PARAM => { }
        _ => { }
    }
}'
$ws.Range("E755").Value = 'fn foo(a: usize, b: usize) -> usize { a + b }
const FOO: fn(usize, usize) -> usize = foo;
fn main() {
    match FOO {
        FOO => {},
        _ => {},
    }
}'
$ws.Range("E756").Value = '
struct NoDerive(u32);
impl PartialEq for NoDerive { fn eq(&self, _: &Self) -> bool { false } }
impl Eq for NoDerive { }
fn main() {
    const INDEX: Option<NoDerive> = [None, Some(NoDerive(10))][0];
    match None { Some(_) => panic!("whoops"), INDEX => dbg!(INDEX), };
}'
$ws.Range("E758").Value = 'extern crate test;

fn name(b: &mut test::Bencher) {
    b.iter(|| 123)
}'
$ws.Range("E759").Value = '
fn x() {}

fn main() {
    x()
}'
$ws.Range("E760").Value = 'use std::arch::asm;

fn main() {
    #[cfg(target_arch="x86_64")]
    unsafe {
        asm!("mov {0}, {0}", in(reg) 0i16);
    }
}

This will produce:
warning: formatting may not be suitable for sub-register argument
 --> src/main.rs:7:19
  |
7 |         asm!("mov {0}, {0}", in(reg) 0i16);
  |                   ^^^  ^^^           ---- for this argument
  |
  = note: `#[warn(asm_sub_register)]` on by default
  = help: use the `x` modifier to have the register formatted as `ax`
  = help: or use the `r` modifier to keep the default formatting of `rax`'
$ws.Range("E761").Value = 'use std::arch::asm;

fn main() {
    #[cfg(target_arch="x86_64")]
    unsafe {
        asm!(
            ".att_syntax",
            "movq %{0}, %{0}", in(reg) 0usize
        );
    }
}

This will produce:
warning: avoid using `.att_syntax`, prefer using `options(att_syntax)` instead
 --> src/main.rs:8:14
  |
8 |             ".att_syntax",
  |              ^^^^^^^^^^^
  |
  = note: `#[warn(bad_asm_style)]` on by default'
$ws.Range("E762").Value = '
unsafe fn foo() {}

unsafe fn bar() {
    foo();
}

fn main() {}'
$ws.Range("E763").Value = 'enum E {
    A,
}

impl Drop for E {
    fn drop(&mut self) {
        println!("Drop");
    }
}

fn main() {
    let e = E::A;
    let i = e as u32;
}'
$ws.Range("E764").Value = '
fn main() {
    let _dangling = 16_usize as *const u8;
}'
$ws.Range("E765").Value = '
fn main() {
    let x: u8 = 37;
    let _addr: usize = &x as *const u8 as usize;
}'
$ws.Range("E769").Value = 'struct X;

impl Default for X {
    fn default() -> Self {
        X
    }
}'
$ws.Range("E770").Value = '
use std::arch::asm;

pub fn default_abi() -> u32 {
    unsafe { asm!("", options(noreturn)); }
}

pub extern "Rust" fn rust_abi() -> u32 {
    unsafe { asm!("", options(noreturn)); }
}'
$ws.Range("E771").Value = '
struct S {}

impl Copy for S {}'
$ws.Range("E772").Value = 'macro_rules! foo {
    () => { true; }
}

fn main() {
    let val = match true {
        true => false,
        _ => foo!()
    };
}'
$ws.Range("E773").Value = 'struct S { /* fields */ }

produces:
warning: derive helper attribute is used before it is introduced
  --> $DIR/legacy-derive-helpers.rs:1:3
   |
 1 | #[serde(rename_all = "camelCase")]
   |   ^^^^^
...
 2 | #[derive(Deserialize)]
   |          ----------- the attribute is introduced here'
$ws.Range("E775").Value = ''
$ws.Range("E776").Value = 'rustc --check-cfg ''names()''

fn foo() {}

This will produce:
warning: unknown condition name used
 --> lint_example.rs:1:7
  |
1 | #[cfg(widnows)]
  |       ^^^^^^^
  |
  = note: `#[warn(unexpected_cfgs)]` on by default'
$ws.Range("E777").Value = 'use foo::NonExhaustiveZst;

struct Bar(u32, ([u32; 0], NonExhaustiveZst));

This will produce:
error: zero-sized fields in repr(transparent) cannot contain external non-exhaustive types
 --> src/main.rs:5:28
  |
5 | struct Bar(u32, ([u32; 0], NonExhaustiveZst));
  |                            ^^^^^^^^^^^^^^^^
  |
note: the lint level is defined here
 --> src/main.rs:1:9
  |
1 | #![deny(repr_transparent_external_private_fields)]
  |         ^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^
  = warning: this was previously accepted by the compiler but is being phased out; it will become a hard error in a future release!
  = note: for more information, see issue #78586 <https://github.com/rust-lang/rust/issues/78586>
  = note: this struct contains `NonExhaustiveZst`, which is marked with `#[non_exhaustive]`, and makes it not a breaking change to become non-zero-sized in the future.'
$ws.Range("E778").Value = 'macro foo() {}'
$ws.Range("E779").Value = 'pub mod foo {
    pub type X = u8;
}

pub mod bar {
    pub type Y = u8;
    pub type X = u8;
}

pub use foo::*;
pub use bar::*;


pub fn main() {}'
$ws.Range("E781").Value = 'use std::thread;

struct Pointer(*mut i32);
unsafe impl Send for Pointer {}

fn main() {
    let mut f = 10;
    let fptr = Pointer(&mut f as *mut i32);
    thread::spawn(move || unsafe {
        *fptr.0 = 20;
    });
}'
$ws.Range("E782").Value = '
extern fn foo() {}'
$ws.Range("E783").Value = '
pub mod submodule {
    #![doc(test(no_crate_inject))]
}'
$ws.Range("E785").Value = '
macro_rules! match_any {
    ( $expr:expr , $( $( $pat:pat )|+ => $expr_arm:expr ),+ ) => {
        match $expr {
            $(
                $( $pat => $expr_arm, )+
            )+
        }
    };
}

fn main() {
    let result: Result<i64, i32> = Err(42);
    let int: i64 = match_any!(result, Ok(i) | Err(i) => i.into());
    assert_eq!(int, 42);
}'
$ws.Range("E786").Value = '
trait Foo {
    fn try_into(self) -> Result<String, !>;
}

impl Foo for &str {
    fn try_into(self) -> Result<String, !> {
        Ok(String::from(self))
    }
}

fn main() {
    let x: String = "3".try_into().unwrap();
    //                  ^^^^^^^^
    // This call to try_into matches both Foo::try_into and TryInto::try_into as
    // `TryInto` has been added to the Rust prelude in 2021 edition.
    println!("{x}");
}'
$ws.Range("E787").Value = '
macro_rules! m {
    (z $x:expr) => ();
}

m!(z"hey");'
$ws.Range("E790").Value = '// crate A
pub enum Bar {
    A,
    B, // added variant in non breaking change
}

// in crate B

match Bar::A {
    Bar::A => {},
    #[warn(non_exhaustive_omitted_patterns)]
    _ => {},
}

This will produce:
warning: reachable patterns not covered of non exhaustive enum
   --> $DIR/reachable-patterns.rs:70:9
   |
LL |         _ => {}
   |         ^ pattern `B` not covered
   |
 note: the lint level is defined here
  --> $DIR/reachable-patterns.rs:69:16
   |
LL |         #[warn(non_exhaustive_omitted_patterns)]
   |                ^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^^
   = help: ensure that all possible cases are being handled by adding the suggested match arms
   = note: the matched value is of type `Bar` and the `non_exhaustive_omitted_patterns` attribute was found'
$ws.Range("E791").Value = 'fn main() {
    println!("{:?}"); // ''‮'');
}'
$ws.Range("E792").Value = 'fn foo() {}

This will produce:
warning: duplicated attribute
 --> src/lib.rs:2:1
  |
2 | #[test]
  | ^^^^^^^
  |
  = note: `#[warn(duplicate_macro_attributes)]` on by default'
$ws.Range("E795").Value = ''
$ws.Range("E796").Value = '
extern "C-unwind" {
    fn foo();
}

fn bar() {
    unsafe { foo(); }
    let ptr: unsafe extern "C-unwind" fn() = foo;
    unsafe { ptr(); }
}'
$ws.Range("E797").Value = 'fn main() {
    let _x = 5;
    println!("{}", _x = 1); // Prints 1, will trigger lint

    println!("{}", _x); // Prints 5, no lint emitted
    println!("{_x}", _x = _x); // Prints 5, no lint emitted
}'
$ws.Range("E798").Value = '
trait Trait {
    fn get<''s>(s: &''s str, _: &''static &''static ()) -> &''static str;
}

impl Trait for () {
    fn get<''s>(s: &''s str, _: &''static &''s ()) -> &''static str {
        s
    }
}

let val = <() as Trait>::get(&String::from("blah blah blah"), &&());
println!("{}", val);'
$ws.Range("E799").Value = 'struct FlexZeroSlice {
    width: u8,
    data: [u8],
}'
$ws.Range("E800").Value = '
macro_rules! myMacro {
   () => {
        // [...]
   }
}
'
